$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.805.57'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '2.469.88'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.09'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.19'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.550'
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").Value = '  +3.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.80'
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +8.49%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '2.848.45'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.73'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '2.462.40'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.780'
$ws.Range("E17").Value = '  +3.57%  '
$ws.Range("D18").Value = '41.725.23'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("E19").Value = '  +2.94%  '
$ws.Range("D20").Value = '0.0₃0951'
$ws.Range("E20").Value = '  +3.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.04'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.48'
$ws.Range("E22").Value = '  +2.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.64'
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.65'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.79'
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.87'
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.77'
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.55'
$ws.Range("E32").Value = '  +2.60%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0765'
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.57'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.63'
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("B36").Value = 'ApeXProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.48'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.90'
$ws.Range("E37").Value = '  -1.57%  '
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.81'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '1.975.62'
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.07'
$ws.Range("E44").Value = '  -3.92%  '
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.04'
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").Value = '2.700.24'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.96'
$ws.Range("E49").Value = '  +0.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.29'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.79'
$ws.Range("E51").Value = '  -0.67%  '
